# The edit re-orders the observation records that live in rows 2-8 of the
# "Artfynd" sheet (row 1 is the header row and is left untouched). Every
# column (A:AY) for a given record travels together - nothing in a row
# is individually changed, the rows are simply permuted:
#
#   old row 2 (id 99019545)   -> new row 7
#   old row 3 (id 103672752)  -> new row 8
#   old row 4 (id 103417728)  -> new row 3
#   old row 5 (id 103417724)  -> new row 4
#   old row 6 (id 103417723)  -> new row 5
#   old row 7 (id 103417731)  -> new row 2
#   old row 8 (id 103417726)  -> new row 6
#
# That permutation decomposes into two cycles: (2 7) and (3 8 6 5 4).
# We realise each cycle with a sequence of Range.Cut(Destination) moves,
# using row 100 (well below the used range) as scratch space so a cycle
# never overwrites data it still needs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cycle (2 7): swap rows 2 and 7 via the scratch row.
$ws.Range("A7:AY7").Cut($ws.Range("A100:AY100"))
$ws.Range("A2:AY2").Cut($ws.Range("A7:AY7"))
$ws.Range("A100:AY100").Cut($ws.Range("A2:AY2"))

# Cycle (3 8 6 5 4): rotate rows 4 -> 3, 5 -> 4, 6 -> 5, 8 -> 6, 3 -> 8.
$ws.Range("A4:AY4").Cut($ws.Range("A100:AY100"))
$ws.Range("A5:AY5").Cut($ws.Range("A4:AY4"))
$ws.Range("A6:AY6").Cut($ws.Range("A5:AY5"))
$ws.Range("A8:AY8").Cut($ws.Range("A6:AY6"))
$ws.Range("A3:AY3").Cut($ws.Range("A8:AY8"))
$ws.Range("A100:AY100").Cut($ws.Range("A3:AY3"))
